$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 6.118831537812738
$ws.Range("D2").Value = 9.23425552603333
$ws.Range("E2").Value = 13.70416625564193
$ws.Range("F2").Value = 32.7118824949064
$ws.Range("G2").Value = 3.653902384284682
$ws.Range("I2").Value = 20.14265615807181
$ws.Range("J2").Value = 9.924023894408338
$ws.Range("K2").Value = 13.07806819628019
$ws.Range("O2").Value = 24.61753418735297

$ws.Range("B3").Value = 5.963427134343079
$ws.Range("D3").Value = 9.170319485480816
$ws.Range("E3").Value = 13.63905738314455
$ws.Range("F3").Value = 32.7755648742236
$ws.Range("G3").Value = 3.656048007801576
$ws.Range("I3").Value = 20.26886351966719
$ws.Range("J3").Value = 9.931410198249292
$ws.Range("K3").Value = 12.64952747064381
$ws.Range("O3").Value = 24.71224037871714

$ws.Range("B4").Value = 5.866506434745174
$ws.Range("D4").Value = 9.132266314600624
$ws.Range("E4").Value = 13.60173273520734
$ws.Range("F4").Value = 32.82375679531887
$ws.Range("G4").Value = 3.657435015174226
$ws.Range("I4").Value = 20.35034102370425
$ws.Range("J4").Value = 9.937570432008439
$ws.Range("K4").Value = 12.37966002121385
$ws.Range("O4").Value = 24.77661982066844

$ws.Range("B5").Value = 5.826688235225682
$ws.Range("D5").Value = 9.117074653519513
$ws.Range("E5").Value = 13.5872014147616
$ws.Range("F5").Value = 32.84567430250678
$ws.Range("G5").Value = 3.658017786141091
$ws.Range("I5").Value = 20.38454883317966
$ws.Range("J5").Value = 9.940489595169689
$ws.Range("K5").Value = 12.2681464675158
$ws.Range("O5").Value = 24.80441632317567

$ws.Range("B6").Value = 5.820058668959794
$ws.Range("D6").Value = 9.114571492559245
$ws.Range("E6").Value = 13.58482982320144
$ws.Range("F6").Value = 32.84945110662441
$ws.Range("G6").Value = 3.658115616702612
$ws.Range("I6").Value = 20.3902898009232
$ws.Range("J6").Value = 9.940999015138296
$ws.Range("K6").Value = 12.24954148333583
$ws.Range("O6").Value = 24.80912606215903

$ws.Range("B7").Value = 5.865970661870869
$ws.Range("D7").Value = 9.132060141560402
$ws.Range("E7").Value = 13.60153399748907
$ws.Range("F7").Value = 32.82404316521678
$ws.Range("G7").Value = 3.657442803476379
$ws.Range("I7").Value = 20.35079828894634
$ws.Range("J7").Value = 9.937608145399492
$ws.Range("K7").Value = 12.37816213278513
$ws.Range("O7").Value = 24.77698837899966

$ws.Range("B8").Value = 6.06559209219463
$ws.Range("D8").Value = 9.211968235785067
$ws.Range("E8").Value = 13.68117354424795
$ws.Range("F8").Value = 32.73194943819577
$ws.Range("G8").Value = 3.654627786661323
$ws.Range("I8").Value = 20.18534693898102
$ws.Range("J8").Value = 9.926233567558622
$ws.Range("K8").Value = 12.93179533244191
$ws.Range("O8").Value = 24.64889293539322

$ws.Range("B9").Value = 6.44285750013217
$ws.Range("D9").Value = 9.377644497397455
$ws.Range("E9").Value = 13.85783390526089
$ws.Range("F9").Value = 32.6237620485286
$ws.Range("G9").Value = 3.649657110087038
$ws.Range("I9").Value = 19.89239192791281
$ws.Range("J9").Value = 9.916809900200848
$ws.Range("K9").Value = 13.95785667248827
$ws.Range("O9").Value = 24.44734872922013

$ws.Range("B10").Value = 6.708570815089224
$ws.Range("D10").Value = 9.504042758453933
$ws.Range("E10").Value = 13.99931269758982
$ws.Range("F10").Value = 32.58874382072103
$ws.Range("G10").Value = 3.646336582321922
$ws.Range("I10").Value = 19.69617512186268
$ws.Range("J10").Value = 9.917718438226348
$ws.Range("K10").Value = 14.66810387590947
$ws.Range("O10").Value = 24.32983930236351

$ws.Range("B11").Value = 6.826436198147164
$ws.Range("D11").Value = 9.56238392400606
$ws.Range("E11").Value = 14.06602591469744
$ws.Range("F11").Value = 32.58251824237121
$ws.Range("G11").Value = 3.644897188565073
$ws.Range("I11").Value = 19.61100173081508
$ws.Range("J11").Value = 9.919825584772003
$ws.Range("K11").Value = 14.98046594786861
$ws.Range("O11").Value = 24.2830786770798

$ws.Range("B12").Value = 6.870597011976444
$ws.Range("D12").Value = 9.584582131748057
$ws.Range("E12").Value = 14.09161019837074
$ws.Range("F12").Value = 32.58155854654683
$ws.Range("G12").Value = 3.644362298057509
$ws.Range("I12").Value = 19.57933356341279
$ws.Range("J12").Value = 9.920866262627847
$ws.Range("K12").Value = 15.09711494732061
$ws.Range("O12").Value = 24.26633910776651

$ws.Range("B13").Value = 6.86110778730867
$ws.Range("D13").Value = 9.579796897014889
$ws.Range("E13").Value = 14.08608614628068
$ws.Range("F13").Value = 32.58170304159873
$ws.Range("G13").Value = 3.644477044510149
$ws.Range("I13").Value = 19.5861278848054
$ws.Range("J13").Value = 9.920631352307872
$ws.Range("K13").Value = 15.07206660403763
$ws.Range("O13").Value = 24.26990115996956

$ws.Range("B14").Value = 6.830079013811375
$ws.Range("D14").Value = 9.564208154651054
$ws.Range("E14").Value = 14.06812440660086
$ws.Range("F14").Value = 32.58241126432225
$ws.Range("G14").Value = 3.64485297914213
$ws.Range("I14").Value = 19.60838465862156
$ws.Range("J14").Value = 9.919906342138662
$ws.Range("K14").Value = 14.99009600449262
$ws.Range("O14").Value = 24.28168208165763

$ws.Range("B15").Value = 6.811010378075533
$ws.Range("D15").Value = 9.554672902758522
$ws.Range("E15").Value = 14.05716366776998
$ws.Range("F15").Value = 32.58302715558472
$ws.Range("G15").Value = 3.645084573727338
$ws.Range("I15").Value = 19.62209370495508
$ws.Range("J15").Value = 9.919493838614736
$ws.Range("K15").Value = 14.93967106798275
$ws.Range("O15").Value = 24.28902439783593

$ws.Range("B16").Value = 6.700804026119915
$ws.Range("D16").Value = 9.500245655117343
$ws.Range("E16").Value = 13.99499875300061
$ws.Range("F16").Value = 32.58934617682368
$ws.Range("G16").Value = 3.646432076982324
$ws.Range("I16").Value = 19.70182342171914
$ws.Range("J16").Value = 9.917614748114833
$ws.Range("K16").Value = 14.6474661613386
$ws.Range("O16").Value = 24.33303036578985

$ws.Range("B17").Value = 6.63239608201459
$ws.Range("D17").Value = 9.46706103545187
$ws.Range("E17").Value = 13.95745351745174
$ws.Range("F17").Value = 32.59571021054875
$ws.Range("G17").Value = 3.647276907895757
$ws.Range("I17").Value = 19.75177990964927
$ws.Range("J17").Value = 9.916895275080545
$ws.Range("K17").Value = 14.46538974312904
$ws.Range("O17").Value = 24.36174492305637

$ws.Range("B18").Value = 6.592768568931913
$ws.Range("D18").Value = 9.448054301299218
$ws.Range("E18").Value = 13.93608110051558
$ws.Range("F18").Value = 32.60028398697025
$ws.Range("G18").Value = 3.647769530552681
$ws.Range("I18").Value = 19.78089835958332
$ws.Range("J18").Value = 9.916640943768567
$ws.Range("K18").Value = 14.35965876453277
$ws.Range("O18").Value = 24.37889065919096

$ws.Range("B19").Value = 6.579304361851587
$ws.Range("D19").Value = 9.441633185944742
$ws.Range("E19").Value = 13.92888351054201
$ws.Range("F19").Value = 32.60198936997083
$ws.Range("G19").Value = 3.647937476140052
$ws.Range("I19").Value = 19.79082353291264
$ws.Range("J19").Value = 9.91658224727979
$ws.Range("K19").Value = 14.32369038073052
$ws.Range("O19").Value = 24.38480394031832

$ws.Range("B20").Value = 6.639707612050133
$ws.Range("D20").Value = 9.470585407547185
$ws.Range("E20").Value = 13.96142735426091
$ws.Range("F20").Value = 32.59493819946728
$ws.Range("G20").Value = 3.647186281334857
$ws.Range("I20").Value = 19.74642214916504
$ws.Range("J20").Value = 9.916955363116578
$ws.Range("K20").Value = 14.48487685132059
$ws.Range("O20").Value = 24.35862298258521

$ws.Range("B21").Value = 6.839206031645507
$ws.Range("D21").Value = 9.568784201472759
$ws.Range("E21").Value = 14.07339161977321
$ws.Range("F21").Value = 32.58216529331402
$ws.Range("G21").Value = 3.644742282247218
$ws.Range("I21").Value = 19.60183144369618
$ws.Range("J21").Value = 9.920112714332019
$ws.Range("K21").Value = 15.01421783250677
$ws.Range("O21").Value = 24.27819543892498

$ws.Range("B22").Value = 6.966821214036613
$ws.Range("D22").Value = 9.633571532586615
$ws.Range("E22").Value = 14.14843204119341
$ws.Range("F22").Value = 32.58196518353595
$ws.Range("G22").Value = 3.643204281174936
$ws.Range("I22").Value = 19.51074260738073
$ws.Range("J22").Value = 9.92359066818312
$ws.Range("K22").Value = 15.35060304065132
$ws.Range("O22").Value = 24.23127382491205

$ws.Range("B23").Value = 6.898976098336389
$ws.Range("D23").Value = 9.598942728268332
$ws.Range("E23").Value = 14.10821664802681
$ws.Range("F23").Value = 32.58132598623772
$ws.Range("G23").Value = 3.644019732963234
$ws.Range("I23").Value = 19.55904727820774
$ws.Range("J23").Value = 9.921605303349065
$ws.Range("K23").Value = 15.17197131225335
$ws.Range("O23").Value = 24.25579891288568

$ws.Range("B24").Value = 6.636402997220897
$ws.Range("D24").Value = 9.468991814909323
$ws.Range("E24").Value = 13.95963011839696
$ws.Range("F24").Value = 32.5952843756295
$ws.Range("G24").Value = 3.647227232049604
$ws.Range("I24").Value = 19.74884315335813
$ws.Range("J24").Value = 9.916927701040231
$ws.Range("K24").Value = 14.47607000187377
$ws.Range("O24").Value = 24.36003242687948

$ws.Range("B25").Value = 6.342620304461604
$ws.Range("D25").Value = 9.33194451829779
$ws.Range("E25").Value = 13.80792907702858
$ws.Range("F25").Value = 32.64524041991208
$ws.Range("G25").Value = 3.650943348893144
$ws.Range("I25").Value = 19.96829131450808
$ws.Range("J25").Value = 9.91798185167532
$ws.Range("K25").Value = 13.68746149291996
$ws.Range("O25").Value = 24.49652509558863
